$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Input change: Arbeidsgiver 2's aktuell manedsinntekt goes from 14000 to 5000 ---
$ws.Range("D3").Value = 5000

# --- 2. Clear the (now unused) C6:F6 formulas, fikser tilfeller hvor nevner er 0 ---
$ws.Range("C6:F6").ClearContents()

# --- 3. B30: avoid relying on MIN(1, ...) against a possibly-zero denominator;
#        divide by the larger of B29/B25 instead ---
$ws.Range("B30").Formula = "=B29/MAX(B29,B25)"
$ws.Range("B30").NumberFormat = "_-* #,##0.000000_-;\-* #,##0.000000_-;_-* ""-""??????_-;_-@_-"

# --- 4. Give D26/E26 (and E27) explicit number formats matching B30's / a higher-precision
#        variant, so the now-visible intermediate values are formatted consistently ---
$ws.Range("D26").NumberFormat = "_-* #,##0.000000_-;\-* #,##0.000000_-;_-* ""-""??????_-;_-@_-"
$ws.Range("E26,E27").NumberFormat = "_-* #,##0.0000000000_-;\-* #,##0.0000000000_-;_-* ""-""??????????_-;_-@_-"

# --- 5. B48: guard against dividing by zero when B26 is 0 ---
# First, move B48's current (fill+border) look over to C48, which used to be blank-styled.
$ws.Range("B48").Copy()
$ws.Range("C48").PasteSpecial(-4122)
$ws.Range("B48").Formula = "=IF(B26 > 0, MIN(1,B47/B26), 0)"
$ws.Range("B48").NumberFormat = "0.0000000000"

# --- 6. Restore default scroll position / move the active selection to C7 ---
$ws.Range("C7").Select()
